$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Means" ---
$ws1 = $wb.Worksheets.Item("Means")

# Row 9 - Total Cancer Risk (per million)
$ws1.Range("B9").Value = 23
$ws1.Range("C9").Value = 33
$ws1.Range("D9").Value = 60
$ws1.Range("E9").Value = 65
$ws1.Range("F9").Value = 66
$ws1.Range("G9").Value = 54

# Row 10 - Total Respiratory (hazard quotient)
$ws1.Range("B10").Value = 0.27
$ws1.Range("C10").Value = 0.37
$ws1.Range("D10").Value = 0.5
$ws1.Range("E10").Value = 0.45
$ws1.Range("F10").Value = 0.44
$ws1.Range("G10").Value = 0.45

# --- Sheet 2: "Standard Deviations" ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# Row 9 - Total Cancer Risk (per million) SD
$ws2.Range("B9").Value = 7.2
$ws2.Range("C9").Value = 9.5
$ws2.Range("E9").Value = 7.4
$ws2.Range("F9").Value = 7.5
$ws2.Range("G9").Value = 11

# Row 10 - Total Respiratory (hazard quotient) SD
$ws2.Range("B10").Value = 0.094
$ws2.Range("C10").Value = 0.083
$ws2.Range("G10").Value = 0.075

$wb.Save()
